# Add an "is_contact" column (H) to the contributors header row, with:
#   - a header comment explaining the column (mirrors the existing header comments)
#   - a boolean (TRUE/FALSE) list data validation applied to the column body
#
# Matches commit "Add is_contact; Fix #604".

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# New header cell, value comes from (new) shared string "is_contact".
$ws.Range("H1").Value = "is_contact"

# Reuse the bold/centered header style already used by A1:G1 (style index 1)
# instead of fabricating a brand-new style entry.
$ws.Range("G1").Copy()
$ws.Range("H1").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# Header comment, same look (Tahoma 8pt) as the other header comments.
$ws.Range("H1").AddComment("Is this individual a contact for DOI purposes?")

# Boolean data validation for the whole column body (row 2 down to the last row).
$dvRange = $ws.Range("H2:H1048576")
$dvRange.Validation.Add(3, 1, 1, '"TRUE,FALSE"')
$dvRange.Validation.ErrorTitle = "Not a boolean"
$dvRange.Validation.ErrorMessage = "The values in this column must be `"TRUE`" or `"FALSE`"."
$dvRange.Validation.IgnoreBlank = $true
$dvRange.Validation.ShowInput = $true
$dvRange.Validation.ShowError = $true
